$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{cell='D2'; value='29.078.54'; numeric=$false}
    @{cell='E2'; value='  -0.05%  '; numeric=$false}
    @{cell='D3'; value='1.820.81'; numeric=$false}
    @{cell='E3'; value='  -0.66%  '; numeric=$false}
    @{cell='D4'; value='0.9987'; numeric=$true}
    @{cell='E4'; value='  -0.23%  '; numeric=$false}
    @{cell='D5'; value='241.15'; numeric=$true}
    @{cell='E5'; value='  -0.83%  '; numeric=$false}
    @{cell='D6'; value='0.6150'; numeric=$true}
    @{cell='E6'; value='  -1.99%  '; numeric=$false}
    @{cell='D7'; value='1.001'; numeric=$true}
    @{cell='E7'; value='  -0.07%  '; numeric=$false}
    @{cell='D8'; value='0.07325'; numeric=$true}
    @{cell='E8'; value='  -2.22%  '; numeric=$false}
    @{cell='D9'; value='0.2877'; numeric=$true}
    @{cell='E9'; value='  -1.50%  '; numeric=$false}
    @{cell='D10'; value='22.87'; numeric=$true}
    @{cell='E10'; value='  -1.49%  '; numeric=$false}
    @{cell='D11'; value='0.07656'; numeric=$true}
    @{cell='E11'; value='  -0.35%  '; numeric=$false}
    @{cell='D12'; value='1.821.51'; numeric=$false}
    @{cell='E12'; value='  -0.63%  '; numeric=$false}
    @{cell='D13'; value='4.940'; numeric=$true}
    @{cell='E13'; value='  -1.25%  '; numeric=$false}
    @{cell='D14'; value='0.6590'; numeric=$true}
    @{cell='E14'; value='  -1.19%  '; numeric=$false}
    @{cell='D15'; value='81.63'; numeric=$true}
    @{cell='E15'; value='  -1.31%  '; numeric=$false}
    @{cell='D16'; value='0.000009021'; numeric=$true}
    @{cell='E16'; value='  -3.84%  '; numeric=$false}
    @{cell='D17'; value='5.823'; numeric=$true}
    @{cell='E17'; value='  -2.56%  '; numeric=$false}
    @{cell='D18'; value='29.058.05'; numeric=$false}
    @{cell='E18'; value='  -0.12%  '; numeric=$false}
    @{cell='D19'; value='2.063.88'; numeric=$false}
    @{cell='E19'; value='  -0.60%  '; numeric=$false}
    @{cell='D20'; value='237.12'; numeric=$true}
    @{cell='E20'; value='  +6.26%  '; numeric=$false}
    @{cell='D21'; value='12.42'; numeric=$true}
    @{cell='E21'; value='  -1.34%  '; numeric=$false}
    @{cell='D22'; value='0.9998'; numeric=$true}
    @{cell='D23'; value='7.105'; numeric=$true}
    @{cell='E23'; value='  +0.20%  '; numeric=$false}
    @{cell='D24'; value='0.9995'; numeric=$true}
    @{cell='E24'; value='  -0.27%  '; numeric=$false}
    @{cell='D25'; value='157.48'; numeric=$true}
    @{cell='E25'; value='  -1.53%  '; numeric=$false}
    @{cell='D26'; value='0.1400'; numeric=$true}
    @{cell='E26'; value='  +0.59%  '; numeric=$false}
    @{cell='D27'; value='8.426'; numeric=$true}
    @{cell='E27'; value='  -0.69%  '; numeric=$false}
    @{cell='D28'; value='17.56'; numeric=$true}
    @{cell='E28'; value='  -1.84%  '; numeric=$false}
    @{cell='E29'; value='  -0.92%  '; numeric=$false}
    @{cell='D30'; value='0.05554'; numeric=$true}
    @{cell='E30'; value='  -1.56%  '; numeric=$false}
    @{cell='D31'; value='4.090'; numeric=$true}
    @{cell='E31'; value='  +0.27%  '; numeric=$false}
    @{cell='D32'; value='4.087'; numeric=$true}
    @{cell='E32'; value='  -1.53%  '; numeric=$false}
    @{cell='D33'; value='1.205'; numeric=$true}
    @{cell='E33'; value='  -0.41%  '; numeric=$false}
    @{cell='B34'; value='ImmutableX'; numeric=$false}
    @{cell='C34'; value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; numeric=$false}
    @{cell='D34'; value='0.7336'; numeric=$true}
    @{cell='E34'; value='  -1.01%  '; numeric=$false}
    @{cell='B35'; value='LidoDAOToken'; numeric=$false}
    @{cell='C35'; value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; numeric=$false}
    @{cell='D35'; value='1.809'; numeric=$true}
    @{cell='E35'; value='  -1.51%  '; numeric=$false}
    @{cell='D36'; value='1.130'; numeric=$true}
    @{cell='E36'; value='  -0.90%  '; numeric=$false}
    @{cell='D37'; value='2.612'; numeric=$true}
    @{cell='E37'; value='  -2.25%  '; numeric=$false}
    @{cell='D38'; value='2.828'; numeric=$true}
    @{cell='D39'; value='1.207.75'; numeric=$false}
    @{cell='E39'; value='  -0.99%  '; numeric=$false}
    @{cell='D40'; value='0.01751'; numeric=$true}
    @{cell='E40'; value='  -1.47%  '; numeric=$false}
    @{cell='D41'; value='6.338'; numeric=$true}
    @{cell='E41'; value='  -2.99%  '; numeric=$false}
    @{cell='D42'; value='0.8912'; numeric=$true}
    @{cell='E42'; value='  +0.23%  '; numeric=$false}
    @{cell='D43'; value='1.001'; numeric=$true}
    @{cell='E43'; value='  -0.09%  '; numeric=$false}
    @{cell='D44'; value='100.93'; numeric=$true}
    @{cell='E44'; value='  -1.03%  '; numeric=$false}
    @{cell='D45'; value='1.968.37'; numeric=$false}
    @{cell='E45'; value='  -0.68%  '; numeric=$false}
    @{cell='D46'; value='64.38'; numeric=$true}
    @{cell='E46'; value='  -1.94%  '; numeric=$false}
    @{cell='D47'; value='0.5079'; numeric=$true}
    @{cell='E47'; value='  -0.28%  '; numeric=$false}
    @{cell='D48'; value='0.00000000117'; numeric=$true}
    @{cell='E48'; value='  -5.02%  '; numeric=$false}
    @{cell='D49'; value='0.3992'; numeric=$true}
    @{cell='E49'; value='  -1.85%  '; numeric=$false}
    @{cell='D50'; value='8.997'; numeric=$true}
    @{cell='E50'; value='  +0.36%  '; numeric=$false}
    @{cell='E51'; value='  -1.29%  '; numeric=$false}
)

foreach ($item in $changes) {
    $rng = $ws.Range($item.cell)
    if ($item.numeric) {
        $rng.Value = "'" + $item.value
    } else {
        $rng.Value = $item.value
    }
    $rng.Style = "Normal"
}
